$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 73 - this shifts the existing rows 73:100 down to 74:101
$ws.Rows.Item(73).Insert()

# Populate the newly inserted row 73 with a new weekly price entry.
# Columns A, B, C, E, F, G, H, I, O, R are constant for this market/product block,
# so reuse the same values as the surrounding rows.
$ws.Cells.Item(73, 1).Value = 2
$ws.Cells.Item(73, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(73, 3).Value = "Coquimbo"
$ws.Cells.Item(73, 4).Value = 44924
$ws.Cells.Item(73, 5).Value = 4
$ws.Cells.Item(73, 6).Value = 100112030
$ws.Cells.Item(73, 7).Value = "Poroto granado"
$ws.Cells.Item(73, 8).Value = "Sin especificar"
$ws.Cells.Item(73, 9).Value = "Primera"
$ws.Cells.Item(73, 10).Value = 760
$ws.Cells.Item(73, 11).Value = 14000
$ws.Cells.Item(73, 12).Value = 15000
$ws.Cells.Item(73, 13).Value = 14500
$ws.Cells.Item(73, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(73, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(73, 16).Value = 967
$ws.Cells.Item(73, 17).Value = 15
$ws.Cells.Item(73, 18).Value = "Hortaliza"
